$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "oom" and "timeout" result categories were merged into a single
# "non-terminated" category after a bug fix. Update every data cell that
# previously read "oom" or "timeout" to read "non-terminated" instead.
$ws.Range("D3").Value = "non-terminated"
$ws.Range("D4").Value = "non-terminated"
$ws.Range("E4").Value = "non-terminated"
$ws.Range("D5").Value = "non-terminated"
$ws.Range("E5").Value = "non-terminated"
$ws.Range("F5").Value = "non-terminated"
$ws.Range("D6").Value = "non-terminated"
$ws.Range("E6").Value = "non-terminated"
$ws.Range("F6").Value = "non-terminated"
$ws.Range("G6").Value = "non-terminated"

# Row 9's label and COUNTIF formulas move from "oom" to the merged
# "non-terminated" category (set the whole row in one go so the shared
# formula keeps a single ref/si group). Borrow the plain data-cell format
# (e.g. from D3) for the label cell, matching the other count-row labels.
$ws.Range("D3").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C9").Value = "non-terminated"
$ws.Range("D9:H9").Formula = '=COUNTIF(D$2:D$6, "non-terminated")'

# Row 10 ("timeout") is no longer its own category - clear its contents
# (label, counts and the row total) while keeping the row itself and the
# D:H count cells' number styling intact.
$ws.Range("C10:H10").ClearContents()

# Re-establish the running total formula for column I over the now-shorter
# I7:I9 range (row 10's total cell is removed entirely below).
$ws.Range("I7:I9").Formula = '=SUM(D7:H7)'

# Fully remove the now-unused I10 total cell (ClearContents would keep an
# empty, styled cell behind; Clear removes it altogether).
$ws.Range("I10").Clear()

# C10 no longer highlights a "timeout" category - match the plain
# (unfilled, bold) label style used elsewhere, e.g. C1, instead of the
# highlighted style it had before.
$ws.Range("C1").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
